# Update LUTs for more precise binning boarders
# The N-column helper formulas subtract 1 from the rounded value (J column)
# before it is emitted into the generated C array text (columns P5 / B16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# N5 is a standalone (non-shared) formula.
$ws.Range("N5").Formula = '=IF(J5<>"",J5-1&", "&IF($H5<>$H6,CHAR(10),""),256^$C$8-1&CHAR(10))'

# N6:N69 share one formula (relative references auto-adjust per row).
$ws.Range("N6:N69").Formula = '=IF(J6<>"",J6-1&", "&IF($H6<>$H7,CHAR(10),""),256^$C$8-1&CHAR(10))'

# N70:N79 share another formula group.
$ws.Range("N70:N79").Formula = '=IF(J70<>"",J70-1&", "&IF($H70<>$H71,CHAR(10),""),256^$C$8-1&CHAR(10))'

# Recalculate so the dependent cells (P5, B16, etc.) pick up the new values.
$excel.Calculate()

# Restore the view/selection state recorded in the edited workbook.
$ws.Activate()
$ws.Range("B16:C16").Select() | Out-Null
